$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Date" column header
$ws.Range("C1").Value = "Date"

# Add date values for the remaining two data rows (as text, matching shared string "1/08/2018")
$ws.Range("C2:C3").NumberFormat = "@"
$ws.Range("C2").Value = "1/08/2018"
$ws.Range("C3").Value = "1/08/2018"
$ws.Range("C2:C3").ClearFormats()

# Remove rows 4 and 5 (trailing data no longer needed)
$ws.Rows("4:5").Delete()

# Update the active selection
$ws.Range("E5").Select() | Out-Null
